$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("installs") holds values like "100,000,000+" stored as text.
# Replace them with plain numeric values (matching column E "minInstalls").
$values = @{
    2  = 100000000
    3  = 1000000000
    4  = 100000000
    5  = 100000000
    6  = 50000000
    7  = 100000000
    8  = 10000000
    9  = 5000000
    10 = 1000000
    11 = 1000000
    12 = 1000000
    13 = 500000000
    14 = 10000000
    15 = 10000000
    16 = 1000000
    17 = 10000000
    18 = 10000000
    19 = 10000000
    20 = 5000000
    21 = 10000000
    22 = 10000000
    23 = 5000000
    24 = 10000000
    25 = 5000000
    26 = 5000000
    27 = 10000000
    28 = 10000000
    29 = 50000000
    30 = 10000000
    31 = 10000000
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}
